$wb = $excel.ActiveWorkbook

# GMXB sheet: insert a new column N for the mortality scalar id
$wsGmxb = $wb.Worksheets.Item("GMXB")
$mortTableColWidth = $wsGmxb.Columns("M:M").ColumnWidth
$wsGmxb.Columns("N:N").Insert()
$wsGmxb.Columns("N:N").ColumnWidth = $mortTableColWidth

$wsGmxb.Range("N1").Value = "mort_scalar_id"
$wsGmxb.Range("N2:N5").Value = "M001"

$wsGmxb.Range("N8").Select() | Out-Null

# ConstParams sheet: selection moved
$wsConst = $wb.Worksheets.Item("ConstParams")
$wsConst.Range("C26").Select() | Out-Null

# RunParams sheet: final active sheet/selection
$wsRun = $wb.Worksheets.Item("RunParams")
$wsRun.Activate() | Out-Null
$wsRun.Range("D13").Select() | Out-Null
